# Auto-generated edit script applying the crypto price/volume refresh diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $range = $ws.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $value
}

$ws.Range("D2").Value = "29.427.13"
$ws.Range("E2").Value = "  +0.11%  "

$ws.Range("D3").Value = "1.849.41"
$ws.Range("E3").Value = "  +0.19%  "

Set-TextValue "D4" "1.000"
$ws.Range("E4").Value = "  +0.06%  "

Set-TextValue "D5" "240.62"
$ws.Range("E5").Value = "  +0.62%  "

Set-TextValue "D6" "0.6278"
$ws.Range("E6").Value = "  -0.71%  "

Set-TextValue "D8" "0.07684"
$ws.Range("E8").Value = "  +1.61%  "

Set-TextValue "D9" "0.2922"
$ws.Range("E9").Value = "  -0.43%  "

Set-TextValue "D10" "24.91"
$ws.Range("E10").Value = "  +1.40%  "

Set-TextValue "D11" "0.07753"

$ws.Range("D12").Value = "1.858.06"
$ws.Range("E12").Value = "  +0.86%  "

Set-TextValue "D13" "5.037"
$ws.Range("E13").Value = "  +0.66%  "

Set-TextValue "D14" "0.00001080"
$ws.Range("E14").Value = "  +3.45%  "

Set-TextValue "D15" "0.6817"
$ws.Range("E15").Value = "  +0.25%  "

Set-TextValue "D16" "83.60"
$ws.Range("E16").Value = "  +0.26%  "

$ws.Range("D17").Value = "2.109.19"
$ws.Range("E17").Value = "  +0.99%  "

Set-TextValue "D18" "6.216"
$ws.Range("E18").Value = "  +0.67%  "

$ws.Range("D19").Value = "29.457.88"
$ws.Range("E19").Value = "  +0.14%  "

Set-TextValue "D20" "228.84"
$ws.Range("E20").Value = "  -0.10%  "

$ws.Range("E21").Value = "  -0.32%  "

$ws.Range("E22").Value = "  -0.02%  "

Set-TextValue "D23" "7.459"
$ws.Range("E23").Value = "  -0.24%  "

Set-TextValue "D25" "157.62"
$ws.Range("E25").Value = "  +0.45%  "

Set-TextValue "D26" "0.1380"
$ws.Range("E26").Value = "  -0.93%  "

Set-TextValue "D27" "8.419"
$ws.Range("E27").Value = "  +0.80%  "

Set-TextValue "D28" "17.75"
$ws.Range("E28").Value = "  +0.86%  "

Set-TextValue "D29" "1.366"
$ws.Range("E29").Value = "  +5.10%  "

Set-TextValue "D30" "1.462"
$ws.Range("E30").Value = "  +0.30%  "

Set-TextValue "D31" "0.05632"
$ws.Range("E31").Value = "  -0.39%  "

Set-TextValue "D32" "4.127"
$ws.Range("E32").Value = "  +0.67%  "

Set-TextValue "D33" "4.049"
$ws.Range("E33").Value = "  +0.62%  "

Set-TextValue "D34" "1.845"

Set-TextValue "D35" "1.163"
$ws.Range("E35").Value = "  +0.47%  "

Set-TextValue "D36" "0.7063"
$ws.Range("E36").Value = "  -0.56%  "

$ws.Range("E37").Value = "  +0.20%  "

$ws.Range("D38").Value = "1.224.84"
$ws.Range("E38").Value = "  -2.01%  "

$ws.Range("E39").Value = "  -0.83%  "

Set-TextValue "D40" "2.759"
$ws.Range("E40").Value = "  -0.24%  "

Set-TextValue "D41" "6.459"
$ws.Range("E41").Value = "  +1.18%  "

Set-TextValue "D42" "0.9030"
$ws.Range("E42").Value = "  +0.04%  "

Set-TextValue "D43" "1.001"
$ws.Range("E43").Value = "  +0.04%  "

$ws.Range("D44").Value = "2.015.91"
$ws.Range("E44").Value = "  +0.88%  "

Set-TextValue "D45" "101.97"
$ws.Range("E45").Value = "  +0.19%  "

Set-TextValue "D46" "66.15"
$ws.Range("E46").Value = "  +0.41%  "

$ws.Range("E47").Value = "  +1.34%  "

$ws.Range("E48").Value = "  +0.49%  "

Set-TextValue "D49" "0.4021"
$ws.Range("E49").Value = "  +0.57%  "

# Rows 50/51: Algorand and EnergySwap swapped rank order
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue "D50" "0.1157"
$ws.Range("E50").Value = "  +3.00%  "

$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D51" "9.021"
$ws.Range("E51").Value = "  +0.46%  "

